$d = $word.ActiveDocument

# The "Delete-" bullet paragraph that lists entities to delete by Id.
# Find the "Loan" word inside that specific paragraph (paragraph 9),
# scoping the search so we do not touch any other "Loan" occurrence
# elsewhere in the document.
$p = $d.Paragraphs.Item(9)
$r = $p.Range.Duplicate

$found = $r.Find.Execute("Loan", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Move the hidden "_GoBack" bookmark from around "Transaction" to
    # wrap the newly struck-through "Loan" (mirrors what Word does when
    # the last edit happens at that location).
    $d.Bookmarks.Add("_GoBack", $r)

    # Mark "Loan" as struck-through (deleted), matching the formatting
    # already used for Account / Customer / Transaction in this list.
    $r.Font.StrikeThrough = 1
}
